$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-15 Tuesday", "2023-08-16 Wednesday"),
    @("72÷4=", "86÷7="),
    @("65÷5=", "25÷3="),
    @("35÷9=", "85÷8="),
    @("83÷2=", "51÷3="),
    @("33÷4=", "99÷2="),
    @("30÷7=", "19÷7="),
    @("37÷4=", "70÷2="),
    @("59÷2=", "80÷4="),
    @("42÷5=", "49÷3="),
    @("84÷8=", "51÷4="),
    @("72÷5=", "69÷7="),
    @("30÷5=", "61÷4="),
    @("38÷7=", "57÷8="),
    @("58÷9=", "45÷8="),
    @("49÷4=", "54÷4="),
    @("20÷2=", "47÷4="),
    @("60÷9=", "16÷6="),
    @("99÷6=", "95÷4="),
    @("73÷4=", "55÷4="),
    @("45÷3=", "53÷3="),
    @("50÷2=", "21÷8="),
    @("54÷9=", "22÷5="),
    @("76÷6=", "21÷6="),
    @("76÷2=", "47÷8="),
    @("92÷5=", "79÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
